# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the 8ace7ffb file on the zh-cn and
# de-de report sheets, and the corresponding "Latest HO Xliff Generate Date"
# roll-up on the Overview sheet, to reflect a fresh handoff-report generation.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 7 is the 8ace7ffb-... file; column H is "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-19 16:46:15"

# de-de sheet: row 7 is the 8ace7ffb-... file; column H is "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-19 16:46:19"

# Overview sheet: row 7 is the 8ace7ffb-... file; column G is
# "Latest HO Xliff Generate Date" -- pick up the newest of the two handoffs.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-19 16:46:19"
